# Update cryptos list data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.812.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("E2").Style = "Normal"

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.934.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E3").Style = "Normal"

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E4").Style = "Normal"

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("E5").Style = "Normal"

# Row 6 - USDC
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E6").Style = "Normal"

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4887"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E7").Style = "Normal"

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2959"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E8").Style = "Normal"

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06872"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("E9").Style = "Normal"

# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("E10").Style = "Normal"

# Row 11 - Litecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "105.73"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.57%  "
$ws.Range("E11").Style = "Normal"

# Row 12 - WrappedEther
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.937.54"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("E12").Style = "Normal"

# Row 13 - TRON
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "TRON"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07785"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("E13").Style = "Normal"

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.333"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.20%  "
$ws.Range("E14").Style = "Normal"

# Row 15 - Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7012"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("E15").Style = "Normal"

# Row 16 - BitcoinCash
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "273.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.31%  "
$ws.Range("E16").Style = "Normal"

# Row 17 - WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.825.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("E17").Style = "Normal"

# Row 18 - ShibaInu
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007706"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("E18").Style = "Normal"

# Row 19 - Uniswap
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.636"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.69%  "
$ws.Range("E19").Style = "Normal"

# Row 20 - Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.22%  "
$ws.Range("E20").Style = "Normal"

# Row 21 - Dai
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E21").Style = "Normal"

# Row 22 - BinanceUSD
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E22").Style = "Normal"

# Row 23 - Chainlink
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.523"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("E23").Style = "Normal"

# Row 24 - Cosmos
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.825"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("E24").Style = "Normal"

# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.95%  "
$ws.Range("E25").Style = "Normal"

# Row 26 - EthereumClassic
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.32%  "
$ws.Range("E26").Style = "Normal"

# Row 27 - LidoDAOToken
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.167"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.79%  "
$ws.Range("E27").Style = "Normal"

# Row 28 - Stellar
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1036"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.05%  "
$ws.Range("E28").Style = "Normal"

# Row 29 - Toncoin
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.385"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.43%  "
$ws.Range("E29").Style = "Normal"

# Row 30 - Filecoin
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.558"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("E30").Style = "Normal"

# Row 31 - PancakeSwap
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.550"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.06%  "
$ws.Range("E31").Style = "Normal"

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.395"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("E32").Style = "Normal"

# Row 33 - Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04887"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("E33").Style = "Normal"

# Row 34 - ImmutableX
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7564"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("E34").Style = "Normal"

# Row 35 - ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.148"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.98%  "
$ws.Range("E35").Style = "Normal"

# Row 36 - Frax
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.001"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E36").Style = "Normal"

# Row 37 - HuobiToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.729"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E37").Style = "Normal"

# Row 38 - VeChain
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("E38").Style = "Normal"

# Row 39 - Aave
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "79.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +7.36%  "
$ws.Range("E39").Style = "Normal"

# Row 40 - MXToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.663"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("E40").Style = "Normal"

# Row 41 - FraxShare
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.476"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("E41").Style = "Normal"

# Row 42 - RenderToken
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.069"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.24%  "
$ws.Range("E42").Style = "Normal"

# Row 43 - TrustWalletToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8905"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("E43").Style = "Normal"

# Row 44 - TheSandbox
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4444"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("E44").Style = "Normal"

# Row 45 - Quant
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "108.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("E45").Style = "Normal"

# Row 46 - Aptos
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.894"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.18%  "
$ws.Range("E46").Style = "Normal"

# Row 47 - PaxDollar
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E47").Style = "Normal"

# Row 48 - Maker
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "985.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("E48").Style = "Normal"

# Row 49 - Algorand
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1244"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("E49").Style = "Normal"

# Row 50 - Elrond
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("E50").Style = "Normal"

# Row 51 - EnergySwap
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.187"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.99%  "
$ws.Range("E51").Style = "Normal"
